# Adjust filenames for relative path
#
# - Rename "Sheet1" -> "data"
# - Add a "groomed_file" column (B) to the data sheet, listing the groomed
#   (distance transform) filename that corresponds to each segmentation file
#   in column A.
# - Update the "studio" sheet's saved tool/view state so it reflects the
#   groom step ("groom" / "Groomed") instead of the original data step.

$wb = $excel.ActiveWorkbook

# --- Rename the first sheet from "Sheet1" to "data" ---
$dataSheet = $wb.Worksheets.Item("Sheet1")
$dataSheet.Name = "data"

# --- studio sheet: update saved tool_state/view_state to reflect groom step ---
$studioSheet = $wb.Worksheets.Item("studio")
$studioSheet.Range("B3").Value = "groom"
$studioSheet.Range("B4").Value = "Groomed"

# --- data sheet: add groomed_file column (B) ---
$dataSheet.Range("B1").Value = "groomed_file"
$dataSheet.Range("B2").Value = "ellipsoid_1_DT.nrrd"
$dataSheet.Range("B3").Value = "ellipsoid_2_DT.nrrd"
$dataSheet.Range("B4").Value = "ellipsoid_3_DT.nrrd"
$dataSheet.Range("B5").Value = "ellipsoid_4_DT.nrrd"
$dataSheet.Range("B6").Value = "ellipsoid_5_DT.nrrd"
$dataSheet.Range("B7").Value = "ellipsoid_6_DT.nrrd"
$dataSheet.Range("B8").Value = "ellipsoid_7_DT.nrrd"
$dataSheet.Range("B9").Value = "ellipsoid_8_DT.nrrd"
$dataSheet.Range("B10").Value = "ellipsoid_9_DT.nrrd"
